$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LFR_RST_ReleaseFish_EDI_Query")

# Add the two new data rows (row 1 is the header row already present)
$ws.Range("A2").Value = 12
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 263

$ws.Range("A3").Value = 12
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 262

# Expand the query's defined name range to cover the newly added rows
$name = $wb.Names.Item("LFR_RST_ReleaseFish_EDI_Query")
$name.RefersTo = "='LFR_RST_ReleaseFish_EDI_Query'!`$A`$1:`$D`$3"
